$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 164, shifting existing rows 164:272 down to 165:273.
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with the new weekly record.
$ws.Cells.Item(164, 1).Value = 5
$ws.Cells.Item(164, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(164, 3).Value = "Maule"
$ws.Cells.Item(164, 4).Value = 44777
$ws.Cells.Item(164, 5).Value = 7
$ws.Cells.Item(164, 6).Value = "Fruta"
$ws.Cells.Item(164, 7).Value = 100108
$ws.Cells.Item(164, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(164, 9).Value = 100108005
$ws.Cells.Item(164, 10).Value = "Piña"
$ws.Cells.Item(164, 11).Value = "Caramelo"
$ws.Cells.Item(164, 12).Value = "Tercera"
$ws.Cells.Item(164, 13).Value = 540
$ws.Cells.Item(164, 14).Value = 18000
$ws.Cells.Item(164, 15).Value = 18000
$ws.Cells.Item(164, 16).Value = 18000
$ws.Cells.Item(164, 17).Value = "$/caja 16 unidades"
$ws.Cells.Item(164, 18).Value = "Ecuador"
$ws.Cells.Item(164, 19).Value = 1125
$ws.Cells.Item(164, 20).Value = 16
